# Add a new test-result row (row 8) to the X103 sheet, matching the
# format of the existing rows above it.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "Add another row to see data shape"
$ws.Range("B8").Value = 0.0083333333333333332
$ws.Range("B8").NumberFormat = $ws.Range("B7").NumberFormat
$ws.Range("C8").Value = "FAIL"

# Column A was best-fit to the longest "Test Name" string (same width
# as the other two sheets use for their identical column A content).
$ws.Columns.Item(1).ColumnWidth = 30.83

# Leave the selection on the newly-added cell, like the author did.
$ws.Range("C8").Select()
